$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.522.98'

$ws.Range("E2").Value = '  -1.85%  '

$ws.Range("D3").Value = '2.633.42'

$ws.Range("E3").Value = '  -1.52%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.99'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -3.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.61'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.647'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  +5.41%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.123'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -4.69%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.81'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.390'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -2.30%  '

$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.63'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -1.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000187'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -6.05%  '

$ws.Range("D15").Value = '3.105.99'

$ws.Range("E15").Value = '  -1.48%  '

$ws.Range("D16").Value = '64.301.05'

$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").Value = '2.632.15'

$ws.Range("E17").Value = '  -1.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.29'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -3.30%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.67'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -2.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.42'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -1.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '347.33'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.94'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -1.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000114'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.74'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +3.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.39'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -3.70%  '

$ws.Range("E27").Value = '  -1.91%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '568.18'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  +6.59%  '

$ws.Range("E29").Value = '  -2.05%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.08'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -2.58%  '

$ws.Range("B33").Value = 'RenderToken'

$ws.Range("C33").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.70'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +3.08%  '

$ws.Range("B34").Value = 'ImmutableX'

$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.72'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -3.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.31'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").Value = '  -3.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.413'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -2.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.10'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.52%  '

$ws.Range("E39").Value = '  -0.15%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '155.03'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -0.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +4.94%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.51'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -2.79%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.99'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -2.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0600'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -1.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.04'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  +1.24%  '

$ws.Range("E47").Value = '  -0.87%  '

$ws.Range("E48").Value = '  +3.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0251'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -2.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.23'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -3.86%  '

$ws.Range("D51").Value = '0.0₆0239'

$ws.Range("E51").Value = '  -6.02%  '
